# Clear the enrollment values in column I (PSL) for rows 2 through 7.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2:I7").Value = ""
